# Rebuild the player roster list on the active sheet to match the new data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data (header stays the same: Oyuncu Adı | Pozisyon | Takım)
$data = @(
    @("Josh Hart",           "SG,SF,PF",  "New York Knicks"),
    @("Lonzo Ball",          "PG",        "Chicago Bulls"),
    @("Michael Porter Jr.",  "SF,PF",     "Denver Nuggets"),
    @("Kristaps Porzingis",  "PF,C",      "Boston Celtics"),
    @("Chris Boucher",       "PF,C",      "Toronto Raptors"),
    @("Victor Wembanyama",   "C",         "San Antonio Spurs"),
    @("Santi Aldama",        "PF,C",      "Memphis Grizzlies"),
    @("Malik Beasley",       "SG,SF",     "Detroit Pistons"),
    @("Donovan Mitchell",    "PG,SG",     "Cleveland Cavaliers"),
    @("Domantas Sabonis",    "C",         "Sacramento Kings"),
    @("Tim Hardaway Jr.",    "SG,SF",     "Detroit Pistons"),
    @("Alperen Sengün",      "C",         "Houston Rockets"),
    @("Andrew Wiggins",      "SF,PF",     "Golden State Warriors"),
    @("Kelly Oubre Jr.",     "SG,SF",     "Philadelphia 76ers"),
    @("Dyson Daniels",       "PG,SG,SF",  "Atlanta Hawks"),
    @("Cam Thomas",          "SG,SF",     "Brooklyn Nets"),
    @("De'Andre Hunter",     "SF,PF",     "Atlanta Hawks"),
    @("Donte DiVincenzo",    "PG,SG,SF",  "Minnesota Timberwolves")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
